$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 2304
$ws.Range("I3").Value = 2468
$ws.Range("I4").Value = 600
$ws.Range("I5").Value = 218
$ws.Range("I6").Value = 2846
$ws.Range("I7").Value = 8436

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 283
$ws.Range("I8").Value = 529
$ws.Range("I11").Value = 141
$ws.Range("I12").Value = 18
$ws.Range("I15").Value = 110
$ws.Range("I18").Value = 61
$ws.Range("I19").Value = 237
$ws.Range("I20").Value = 223
$ws.Range("I23").Value = 75
$ws.Range("I26").Value = 9
$ws.Range("I27").Value = 84
$ws.Range("I29").Value = 544
$ws.Range("I33").Value = 385
$ws.Range("I34").Value = 33
$ws.Range("I37").Value = 270
$ws.Range("I41").Value = 37
$ws.Range("I42").Value = 292
$ws.Range("I43").Value = 73
$ws.Range("I51").Value = 77
$ws.Range("I52").Value = 168
$ws.Range("I55").Value = 93
$ws.Range("I57").Value = 27
$ws.Range("I60").Value = 44
$ws.Range("I63").Value = 36
$ws.Range("I64").Value = 80
$ws.Range("I65").Value = 199
$ws.Range("I66").Value = 18
$ws.Range("I67").Value = 327
$ws.Range("I79").Value = 220
$ws.Range("I80").Value = 29
$ws.Range("I84").Value = 63
$ws.Range("I85").Value = 391
$ws.Range("I86").Value = 48
$ws.Range("I87").Value = 16
$ws.Range("I90").Value = 95
$ws.Range("I92").Value = 26
$ws.Range("I93").Value = 51
$ws.Range("I95").Value = 142
$ws.Range("I99").Value = 158
$ws.Range("I101").Value = 8436

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 97
$ws.Range("I3").Value = 154
$ws.Range("I7").Value = 391

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 42
$ws.Range("I3").Value = 68
$ws.Range("I7").Value = 168

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 59
$ws.Range("I7").Value = 141

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 143
$ws.Range("I5").Value = 17
$ws.Range("I7").Value = 529

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 85
$ws.Range("I7").Value = 283

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 85
$ws.Range("I3").Value = 81
$ws.Range("I7").Value = 270

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 42
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 327

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I6").Value = 65
$ws.Range("I7").Value = 199

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I3").Value = 58
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 134
$ws.Range("I4").Value = 22
$ws.Range("I7").Value = 385

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 187
$ws.Range("I7").Value = 544

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I6").Value = 68
$ws.Range("I7").Value = 237

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 25
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I2").Value = 14
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I6").Value = 77
$ws.Range("I7").Value = 292

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 220

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 80

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 223

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I6").Value = 11
$ws.Range("I7").Value = 33

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 9

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 18

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I3").Value = 4
$ws.Range("I7").Value = 26

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I6").Value = 7
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I2").Value = 8
$ws.Range("I7").Value = 27

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 73

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I3").Value = 7
$ws.Range("I7").Value = 29

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 18

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I6").Value = 7
$ws.Range("I7").Value = 16
